# Updated cryptos list on Mon Sep 23 03:38:30 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for every coin row,
# and re-sorts a handful of rows whose rank order changed (B/C/D/E all
# updated together for those rows: Fetch.AI/SuiNetwork swap rows 24-25,
# InternetComputer(DFINITY)/Aptos swap rows 27-28,
# PolygonEcosystemToken/FirstDigitalUSD swap rows 35-36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (not an auto-coerced number/date).
# Many "Price" values look numeric ("1.00", "598.66", ...) and a plain
# `.Value = ...` assignment would make Excel store them as numbers, which
# would lose the original formatting (e.g. "1.00" -> 1). Forcing the
# cell to Text format for the duration of the write keeps it a string,
# then resetting the style back to Normal keeps formatting identical to
# the surrounding (untouched) cells.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '64.624.19'
Set-TextValue 'E2' '  +2.36%  '

# Row 3
Set-TextValue 'D3' '2.682.59'
Set-TextValue 'E3' '  +3.02%  '

# Row 4
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  +0.01%  '

# Row 5
Set-TextValue 'D5' '598.66'
Set-TextValue 'E5' '  +2.58%  '

# Row 6
Set-TextValue 'D6' '148.69'
Set-TextValue 'E6' '  +0.38%  '

# Row 7
Set-TextValue 'D7' '1.00'
Set-TextValue 'E7' '  +0.06%  '

# Row 8
Set-TextValue 'D8' '0.595'
Set-TextValue 'E8' '  -0.64%  '

# Row 9
Set-TextValue 'D9' '0.109'
Set-TextValue 'E9' '  +0.66%  '

# Row 10
Set-TextValue 'D10' '5.69'
Set-TextValue 'E10' '  +0.08%  '

# Row 11
Set-TextValue 'E11' '  -0.08%  '

# Row 12
Set-TextValue 'D12' '0.360'
Set-TextValue 'E12' '  +1.50%  '

# Row 13
Set-TextValue 'D13' '28.11'
Set-TextValue 'E13' '  +2.94%  '

# Row 14
Set-TextValue 'D14' '3.159.18'
Set-TextValue 'E14' '  +2.92%  '

# Row 15
Set-TextValue 'D15' '64.472.71'
Set-TextValue 'E15' '  +2.27%  '

# Row 16
Set-TextValue 'E16' '  +0.32%  '

# Row 17
Set-TextValue 'D17' '2.658.43'
Set-TextValue 'E17' '  +2.02%  '

# Row 18
Set-TextValue 'D18' '11.47'
Set-TextValue 'E18' '  +0.77%  '

# Row 19
Set-TextValue 'D19' '348.45'
Set-TextValue 'E19' '  +1.39%  '

# Row 20
Set-TextValue 'D20' '4.43'
Set-TextValue 'E20' '  +0.38%  '

# Row 21
Set-TextValue 'D21' '6.92'
Set-TextValue 'E21' '  +1.97%  '

# Row 22
Set-TextValue 'D22' '1.00'
Set-TextValue 'E22' '  +0.03%  '

# Row 23
Set-TextValue 'D23' '69.37'
Set-TextValue 'E23' '  +3.20%  '

# Row 24
Set-TextValue 'B24' 'Fetch.AI'
Set-TextValue 'C24' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D24' '1.68'
Set-TextValue 'E24' '  +5.23%  '

# Row 25
Set-TextValue 'B25' 'SuiNetwork'
Set-TextValue 'C25' 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue 'D25' '1.62'
Set-TextValue 'E25' '  +10.59%  '

# Row 26
Set-TextValue 'E26' '  -1.38%  '

# Row 27
Set-TextValue 'B27' 'InternetComputer(DFINITY)'
Set-TextValue 'C27' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D27' '8.54'
Set-TextValue 'E27' '  +1.42%  '

# Row 28
Set-TextValue 'B28' 'Aptos'
Set-TextValue 'C28' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D28' '8.01'
Set-TextValue 'E28' '  +1.29%  '

# Row 29
Set-TextValue 'D29' '0.998'
Set-TextValue 'E29' '  -0.18%  '

# Row 30
Set-TextValue 'D30' '535.26'
Set-TextValue 'E30' '  +15.43%  '

# Row 31
Set-TextValue 'E31' '  +3.61%  '

# Row 32
Set-TextValue 'D32' '1.79'
Set-TextValue 'E32' '  +11.42%  '

# Row 33
Set-TextValue 'D33' '0.0₃0830'
Set-TextValue 'E33' '  +0.76%  '

# Row 34
Set-TextValue 'D34' '175.52'
Set-TextValue 'E34' '  -0.77%  '

# Row 35
Set-TextValue 'B35' 'FirstDigitalUSD'
Set-TextValue 'C35' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D35' '1.00'
Set-TextValue 'E35' '  +0.11%  '

# Row 36
Set-TextValue 'B36' 'PolygonEcosystemToken'
Set-TextValue 'C36' 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D36' '0.403'
Set-TextValue 'E36' '  -0.82%  '

# Row 37
Set-TextValue 'D37' '19.37'
Set-TextValue 'E37' '  +0.64%  '

# Row 38
Set-TextValue 'D38' '4.72'
Set-TextValue 'E38' '  +2.53%  '

# Row 39
Set-TextValue 'D39' '1.78'
Set-TextValue 'E39' '  +4.41%  '

# Row 40
Set-TextValue 'D40' '173.00'
Set-TextValue 'E40' '  +8.11%  '

# Row 41
Set-TextValue 'E41' '  -0.11%  '

# Row 42
Set-TextValue 'D42' '40.74'
Set-TextValue 'E42' '  +3.11%  '

# Row 43
Set-TextValue 'D43' '3.81'
Set-TextValue 'E43' '  +0.17%  '

# Row 44
Set-TextValue 'D44' '22.03'
Set-TextValue 'E44' '  +5.12%  '

# Row 45
Set-TextValue 'D45' '0.636'
Set-TextValue 'E45' '  -0.63%  '

# Row 46
Set-TextValue 'D46' '0.0553'
Set-TextValue 'E46' '  +1.09%  '

# Row 47
Set-TextValue 'E47' '  +1.71%  '

# Row 48
Set-TextValue 'D48' '0.0967'
Set-TextValue 'E48' '  -0.85%  '

# Row 49
Set-TextValue 'D49' '19.01'
Set-TextValue 'E49' '  +2.06%  '

# Row 50
Set-TextValue 'D50' '1.79'
Set-TextValue 'E50' '  +2.95%  '

# Row 51
Set-TextValue 'E51' '  -0.45%  '
